$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Iron Workers Local 201 Washington DC" has been resolved/removed -> delete its row.
# This shifts rows 5-7 up to rows 4-6.
$ws.Rows("4:4").Delete()

# Update status notes and fill in newly tracked custodian/status info.
$ws.Range("C4").Value = "Call Jessica Schneider"
$ws.Range("B5").Value = "BNY"
$ws.Range("C5").Value = "Should be on BNY"
$ws.Range("B6").Value = "Regions Bank"
$ws.Range("C3").Value = "Call Admin back"

$ws.Range("C4").Select()

# Portrait page setup (was unset before).
$ws.PageSetup.Orientation = 1
